$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in FEATURES values B21:B27
$ws.Range("B21").Value = 1
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = 0

# Number format changes: progress ratio column F changes from percent to 0.00
$ws.Range("F2:F16").NumberFormat = "0.00"

# B20:B27 changes from General to 0.00
$ws.Range("B20:B27").NumberFormat = "0.00"

# B28 (TOTAL) changes to percent format
$ws.Range("B28").NumberFormat = "0%"

# Update selection
$ws.Range("E26").Select() | Out-Null
